# #5: insurance, claim, debt, investment done
# Rebuild the "債務" (debt) sheet (7th worksheet) so it carries the full
# common schema (property_category/category/date/legislator_name/
# legislator_id/source_file/index) alongside its species-specific columns,
# matching the other detail sheets (land/building/car/cash/deposit/stock).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(7)

# --- Row 1: header labels (B1:G1 previously held stray data values) ---
$ws.Cells.Item(1,2).Value = "species"
$ws.Cells.Item(1,3).Value = "debtor"
$ws.Cells.Item(1,4).Value = "owner"
$ws.Cells.Item(1,5).Value = "total"
$ws.Cells.Item(1,6).Value = "register_date"
$ws.Cells.Item(1,7).Value = "register_reason"
$ws.Cells.Item(1,8).Value = "property_category"
$ws.Cells.Item(1,9).Value = "category"
$ws.Cells.Item(1,10).Value = "date"
$ws.Cells.Item(1,11).Value = "legislator_name"
$ws.Cells.Item(1,12).Value = "legislator_id"
$ws.Cells.Item(1,13).Value = "source_file"
$ws.Cells.Item(1,14).Value = "index"

# --- Row 2: data values for the new trailing columns ---
$ws.Cells.Item(2,8).Value = "debt"
$ws.Cells.Item(2,9).Value = "normal"
$ws.Cells.Item(2,10).Value = "2013-12-24"
$ws.Cells.Item(2,11).Value = "林國正"
$ws.Cells.Item(2,12).Value = 1742
$ws.Cells.Item(2,13).Value = "tmp399c1"
$ws.Cells.Item(2,14).Value = 101

# --- Formatting: extend the existing header/data styles across the new
#     columns instead of rebuilding them property-by-property (which would
#     leave stray unused style entries behind). ---
$ws.Range("B1").Copy()
$ws.Range("H1:N1").PasteSpecial(-4122)
$ws.Range("B2").Copy()
$ws.Range("H2:N2").PasteSpecial(-4122)

$ws.Range("A1").Select()
